$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "What is the most number of tables that I can have in Geo?"
$ws.Range("B14").Value = "llama3.2:latest"
$ws.Range("C14").Value = "The GEO software does not specify a maximum limit for the number of tables that can be put in one ODF file. However, it mentions that there are software limits on this topic."

$ws.Range("A15").Value = "Bullet point list types of limits in GEO."
$ws.Range("B15").Value = "llama3.2:latest"
$ws.Range("C15").Value = "Here are the types of limits in GEO:`n• Modifiers`n• Lithologies`n• Symbols`n• Texts`n• Lines"

$ws.Range("A16").Value = "Bullet point list types of limits in GEO."
$ws.Range("B16").Value = "llama3.2:latest"
$ws.Range("C16").Value = "Here are the types of limits in GEO:`n• Modifiers`n• Lithologies`n• Symbols`n• Texts`n• Lines"
